$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the existing header style (s="1", bold + border + centered) before
# the source cells get overwritten, so the new header row re-uses it
# exactly like the old C3:E3 header row did.
$ws.Range("C3").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

# Now that the format has been copied, remove the old data block entirely.
$ws.Range("C3:E4").Clear()

# New header row (row 1): colour labels + the computed "Hue" column.
$ws.Range("A1").Value = "Blue"
$ws.Range("B1").Value = "Green"
$ws.Range("C1").Value = "Red"
$ws.Range("D1").Value = "Hue"

# New data row (row 2).
$ws.Range("A2").Value = 86.87834324477781
$ws.Range("B2").Value = 118.4321638824971
$ws.Range("C2").Value = 148.4725517736007
$ws.Range("D2").Value = 86.87834324477781
